$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition): refresh "want to go" counters ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 319
$wsExpo.Range("F3").Value = 1328

# --- Sheet "演出" (Performance): the single listed event has passed -
# remove its data row, keeping only the header.
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Rows.Item(2).Delete()

# --- Sheet "全部类型" (All types): the same event that left "演出" also
# needs to disappear here. Its row-index column (A) is a plain sequential
# counter (0 for the header, 1,2,3,... for each data row) that must stay
# sequential, so instead of doing a raw row delete/shift (which would drag
# the old index numbers down with the rest of the row), copy each
# following row's B:I content up into the row above, then drop the now
# duplicated trailing row.
$wsAll = $wb.Worksheets.Item("全部类型")
$lastRow = $wsAll.Cells.Item($wsAll.Rows.Count, 1).End(-4162).Row
$wsAll.Range("B3:I" + $lastRow).Copy()
$wsAll.Range("B2").PasteSpecial(-4163)
$wsAll.Rows.Item($lastRow).Delete()

# Same counter refresh as "展览" for the two events that shifted up.
$wsAll.Range("F2").Value = 319
$wsAll.Range("F3").Value = 1328
